$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Terrain/event data update -----------------------------------------
# The "Event" column (J) for the three data rows previously referenced a
# combined terrain/event id string "1001-1002-1005". The preliminary
# terrain system adds a new terrain node (1009) into that chain, so every
# row now shares the single updated definition "1001-1002-1009-1005".
$newEvent = "1001-1002-1009-1005"
$ws.Range("J3").Value = $newEvent
$ws.Range("J4").Value = $newEvent
$ws.Range("J5").Value = $newEvent

# Give the updated event cells their own explicit font run (11pt 宋体)
# instead of inheriting the implicit default font through style 1.
$eventCells = $ws.Range("J3:J5")
$eventCells.Font.Name = "宋体"
$eventCells.Font.Size = 11

# --- Selection / view state ---------------------------------------------
$ws.Range("J5").Select()
